# Screw Improve & First Complete Chapter
# 1. Chapter 1 steps completed -> add two new rows describing the oil filter
#    wrench steps (data rows, pushed in just after the header/existing
#    "10101"/"10102" rows and before the old "Success0" row).
# 2. A couple of small bugfixes to the sheet layout/selection that came
#    along with the author's resave (column width, selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two new data rows at row 5 (everything from the old row 5
#        downward shifts down by two rows) ------------------------------
$ws.Rows("5:6").Insert()

# New row 5: EN-726-A oil filter wrench placement step
$ws.Range("A5").Value = 10103
$ws.Range("B5").Value = "在机油滤清器处放置EN-726-A机油滤清器扳手"
$ws.Range("C5").Value = 10103

# New row 6: large ratchet wrench removal step
$ws.Range("A6").Value = 10104
$ws.Range("B6").Value = "使用大号棘轮扳手拆卸机油滤清器"
$ws.Range("C6").Value = 10104

# --- 2. Layout / selection bugfixes -------------------------------------
# Column B widened to fit the newly-added (longer) text.
$ws.Columns("B:B").ColumnWidth = 43.43

# Selected cell moved to C7 (the row the old "Success0" entry now sits on),
# and the view no longer pins a frozen/scrolled top-left cell.
$ws.Range("C7").Select()
